$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Korean translation text updates (column F) ---
$ws.Range("F2").Value = "바이오테크 확장 - 코어 -사용 금지-"
$ws.Range("F3").Value = "바이오테크 확장 - 공용"
$ws.Range("F4").Value = "털관련 코어 유전자"
$ws.Range("F5").Value = "바이오테크 확장을 위한 핵심 유전자. 사용하지 마세요!"
$ws.Range("F6").Value = "착유 코어 유전자"
$ws.Range("F7").Value = "바이오테크 확장을 위한 핵심 유전자. 사용하지 마세요!"
$ws.Range("F8").Value = "야수성"
$ws.Range("F9").Value = "이 유전자의 보유자는 타고난 신체를 본능적으로 잘 다루는 전문가입니다. 맨손일 경우 근접 공격력이 증가합니다. 이는 다른 근접 공격력 관련 유전자와 중첩됩니다."
$ws.Range("F11").Value = "이 유전자의 보유자는 식물성 식품을 제대로 소화할 수 없습니다. 섭취할 경우 반드시 식중독에 걸립니다. 하지만 시체를 먹는 것에는 거부감이 적습니다."
$ws.Range("F13").Value = "이 유전자의 보유자는 고기를 제대로 소화할 수 없습니다. 섭취할 경우 반드시 식중독에 걸립니다."
$ws.Range("F14").Value = "대머리 털복숭이"
$ws.Range("F15").Value = "변종 털복숭이"
$ws.Range("F16").Value = "이 유전자의 보유자는 몸 전체에 두꺼운 털이 자라서 추운 온도로부터 보호합니다. 이 변종은 얼굴과 머리 부분의 털이 적습니다."
$ws.Range("F18").Value = "궁수"
$ws.Range("F23").Value = "강궁수"
$ws.Range("F24").Value = "우두머리 광전사"
$ws.Range("F25").Value = "우두머리 궁사"
$ws.Range("F27").Value = "벌목꾼"
$ws.Range("F29").Value = "주민"
$ws.Range("F30").Value = "돌격병"

# --- Column width tweaks ---
# Column F widened (target stored width 40.85; engine snaps width to
# nearest 1/7 character-width pixel grid, so feed it the ColumnWidth that
# lands closest to that stored value).
$ws.Columns("F").ColumnWidth = 40.142857142857146

# Far-right columns XFB:XFC (16382:16383) get an explicit narrower custom
# width, split out from the big default-width block that covers the
# remaining unused columns.
$ws.Columns("XFB:XFC").ColumnWidth = 8.857142857142858

# --- Active selection moved to G30 ---
$ws.Range("G30").Select()
